$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 92, shifting existing rows 92-170 down to 93-171.
$ws.Rows.Item(92).Insert()

# Populate the newly inserted row with the new data record.
$ws.Cells.Item(92, 1).Value = 8
$ws.Cells.Item(92, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(92, 3).Value = "Coquimbo"
$ws.Cells.Item(92, 4).Value = 44447
$ws.Cells.Item(92, 5).Value = 4
$ws.Cells.Item(92, 6).Value = 100112032
$ws.Cells.Item(92, 7).Value = "Zapallo italiano"
$ws.Cells.Item(92, 8).Value = "Sin especificar"
$ws.Cells.Item(92, 9).Value = "Primera"
$ws.Cells.Item(92, 10).Value = 600
$ws.Cells.Item(92, 11).Value = 14000
$ws.Cells.Item(92, 12).Value = 15000
$ws.Cells.Item(92, 13).Value = 14500
$ws.Cells.Item(92, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(92, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(92, 16).Value = 290
$ws.Cells.Item(92, 17).Value = 50
$ws.Cells.Item(92, 18).Value = "Hortaliza"
